# Add the new "game_outcome" worksheet at the end of the workbook
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "game_outcome"

# Fill header row (order chosen so shared-string indices line up: a_team, h_team, a_score, h_score, HOU, MIN, winner)
$newSheet.Range("A1").Value = "id"
$newSheet.Range("B1").Value = "a_team"
$newSheet.Range("C1").Value = "h_team"
$newSheet.Range("D1").Value = "a_score"
$newSheet.Range("E1").Value = "h_score"

# Fill data row values that introduce new shared strings next
$newSheet.Range("B2").Value = "HOU"
$newSheet.Range("C2").Value = "MIN"

# Last new shared string
$newSheet.Range("F1").Value = "winner"

# Remaining (non-string / reused-string) data cells
$newSheet.Range("A2").Value = 1
$newSheet.Range("D2").Value = 3
$newSheet.Range("E2").Value = 5
$newSheet.Range("F2").Value = "home"

# Update sheet selections / active tab to match the target workbook state
$newSheet.Activate()
$newSheet.Range("H6").Select()

$wsOutcome = $wb.Worksheets.Item("outcome")
$wsOutcome.Activate()
$wsOutcome.Range("B1:B3").Select()

$wsBets = $wb.Worksheets.Item("bets")
$wsBets.Activate()
$wsBets.Range("I14").Select()
